$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 103
$ws.Range("H103").Value = 829.6667
$ws.Range("I103").Value = 1309.8
$ws.Range("K103").Value = 3929.4
$ws.Range("M103").Value = -3343.4
# Row 113
$ws.Range("H113").Value = 2400
$ws.Range("I113").Value = 2400
$ws.Range("K113").Value = 2400
$ws.Range("M113").Value = 854
# Row 117
$ws.Range("H117").Value = 249335.5
$ws.Range("J117").Value = 249335.5
$ws.Range("L117").Value = 249335.5
$ws.Range("N117").Value = -258513.5
# Row 137
$ws.Range("H137").Value = 1553.3125
$ws.Range("I137").Value = 1505
$ws.Range("J137").Value = 1659.6
$ws.Range("K137").Value = 4515
$ws.Range("L137").Value = 4978.799999999999
$ws.Range("M137").Value = -1965
$ws.Range("N137").Value = -10078.8

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 6670
$ws.Range("I45").Value = 9197.666999999999
$ws.Range("J45").Value = 2878.5
$ws.Range("K45").Value = 9197.666999999999
$ws.Range("L45").Value = 2878.5
$ws.Range("M45").Value = -8820.666999999999
$ws.Range("N45").Value = -3632.5
# Row 74
$ws.Range("H74").Value = 62501940
$ws.Range("J74").Value = 825
$ws.Range("L74").Value = 825
$ws.Range("N74").Value = -2573
# Row 77
$ws.Range("H77").Value = 62501940
$ws.Range("J77").Value = 825
$ws.Range("L77").Value = 4125
$ws.Range("N77").Value = -12861
# Row 97
$ws.Range("H97").Value = 614.8
$ws.Range("I97").Value = 614.8
$ws.Range("K97").Value = 614.8
$ws.Range("M97").Value = -118.8
# Row 102
$ws.Range("I102").Value = 16675665
$ws.Range("K102").Value = 16675665
$ws.Range("M102").Value = -16674043
# Row 110
$ws.Range("H110").Value = 45982.953
$ws.Range("I110").Value = 50544.3
$ws.Range("J110").Value = 369.5
$ws.Range("K110").Value = 50544.3
$ws.Range("L110").Value = 369.5
$ws.Range("M110").Value = -48499.3
$ws.Range("N110").Value = -4459.5
# Row 122
$ws.Range("H122").Value = 4106.5835
$ws.Range("I122").Value = 2207.182
$ws.Range("K122").Value = 6621.545999999999
$ws.Range("M122").Value = -4171.545999999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 19823.666
$ws.Range("I26").Value = 19823.666
$ws.Range("K26").Value = 19823.666
$ws.Range("M26").Value = -19531.666
# Row 94
$ws.Range("H94").Value = 1455.409
$ws.Range("I94").Value = 1526.95
$ws.Range("K94").Value = 1526.95
$ws.Range("M94").Value = -1075.95
# Row 96
$ws.Range("H96").Value = 56666.668
$ws.Range("I96").Value = 56666.668
$ws.Range("K96").Value = 56666.668
$ws.Range("M96").Value = -53920.668
# Row 107
$ws.Range("H107").Value = 114655.11
$ws.Range("I107").Value = 4092.4285
$ws.Range("K107").Value = 4092.4285
$ws.Range("M107").Value = -2172.4285
# Row 134
$ws.Range("H134").Value = 36430540
$ws.Range("J134").Value = 5000
$ws.Range("L134").Value = 15000
$ws.Range("N134").Value = -20070

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2778.868
$ws.Range("I31").Value = 2579.8538
$ws.Range("J31").Value = 3458.8333
$ws.Range("K31").Value = 2579.8538
$ws.Range("L31").Value = 3458.8333
$ws.Range("M31").Value = -2284.8538
$ws.Range("N31").Value = -4048.8333
# Row 34
$ws.Range("H34").Value = 2778.868
$ws.Range("I34").Value = 2579.8538
$ws.Range("J34").Value = 3458.8333
$ws.Range("K34").Value = 2579.8538
$ws.Range("L34").Value = 3458.8333
$ws.Range("M34").Value = -2377.8538
$ws.Range("N34").Value = -3862.8333
# Row 107
$ws.Range("H107").Value = 865998.8
$ws.Range("I107").Value = 1111499
$ws.Range("K107").Value = 1111499
$ws.Range("M107").Value = -1109579
# Row 122
$ws.Range("H122").Value = 3601.375
$ws.Range("I122").Value = 3468.6667
$ws.Range("K122").Value = 10406.0001
$ws.Range("M122").Value = -7956.000100000001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 246.375
$ws.Range("I12").Value = 148
$ws.Range("J12").Value = 279.16666
$ws.Range("K12").Value = 444
$ws.Range("L12").Value = 837.4999799999999
$ws.Range("M12").Value = -271
$ws.Range("N12").Value = -1183.49998
# Row 68
$ws.Range("H68").Value = 1993.2858
$ws.Range("I68").Value = 1193.75
$ws.Range("K68").Value = 3581.25
$ws.Range("M68").Value = -2770.25
# Row 71
$ws.Range("H71").Value = 1993.2858
$ws.Range("I71").Value = 1193.75
$ws.Range("K71").Value = 10743.75
$ws.Range("M71").Value = -6687.75
# Row 137
$ws.Range("H137").Value = 11113911
$ws.Range("J137").Value = 2966.6667
$ws.Range("L137").Value = 8900.000100000001
$ws.Range("N137").Value = -19100.0001

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 4357.6294
$ws.Range("I122").Value = 2822.4285
$ws.Range("J122").Value = 9730.833000000001
$ws.Range("K122").Value = 8467.2855
$ws.Range("L122").Value = 29192.499
$ws.Range("M122").Value = -6017.2855
$ws.Range("N122").Value = -34092.499

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 36
$ws.Range("H36").Value = 90000
$ws.Range("J36").Value = 90000
$ws.Range("L36").Value = 90000
$ws.Range("N36").Value = -91124
# Row 40
$ws.Range("H40").Value = 2470.8572
$ws.Range("I40").Value = 1474.25
$ws.Range("K40").Value = 1474.25
$ws.Range("M40").Value = -1338.25
# Row 82
$ws.Range("H82").Value = 961
$ws.Range("I82").Value = 990
$ws.Range("J82").Value = 903
$ws.Range("K82").Value = 990
$ws.Range("L82").Value = 903
$ws.Range("M82").Value = -629
$ws.Range("N82").Value = -1625
# Row 85
$ws.Range("H85").Value = 961
$ws.Range("I85").Value = 990
$ws.Range("J85").Value = 903
$ws.Range("K85").Value = 990
$ws.Range("L85").Value = 903
$ws.Range("M85").Value = 258
$ws.Range("N85").Value = -3399
# Row 93
$ws.Range("H93").Value = 1712.4
$ws.Range("I93").Value = 1241.4546
$ws.Range("J93").Value = 5166
$ws.Range("K93").Value = 1241.4546
$ws.Range("L93").Value = 5166
$ws.Range("M93").Value = 6.545399999999972
$ws.Range("N93").Value = -7662
# Row 122
$ws.Range("H122").Value = 11719.2
$ws.Range("I122").Value = 11719.2
$ws.Range("K122").Value = 35157.60000000001
$ws.Range("M122").Value = -32707.60000000001
# Row 132
$ws.Range("H132").Value = 11713593
$ws.Range("I132").Value = 18466950
$ws.Range("K132").Value = 55400850
$ws.Range("M132").Value = -55398320
# Row 136
$ws.Range("H136").Value = 3564.2144
$ws.Range("I136").Value = 3825
$ws.Range("J136").Value = 1999.5
$ws.Range("K136").Value = 11475
$ws.Range("L136").Value = 5998.5
$ws.Range("M136").Value = -8925
$ws.Range("N136").Value = -11098.5

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 80
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 81
$ws.Range("H81").Value = 2799.2727
$ws.Range("I81").Value = 3016
$ws.Range("J81").Value = 2420
$ws.Range("K81").Value = 6032
$ws.Range("L81").Value = 4840
$ws.Range("M81").Value = -4971
$ws.Range("N81").Value = -6962
# Row 83
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 84
$ws.Range("H84").Value = 2799.2727
$ws.Range("I84").Value = 3016
$ws.Range("J84").Value = 2420
$ws.Range("K84").Value = 30160
$ws.Range("L84").Value = 24200
$ws.Range("M84").Value = -24856
$ws.Range("N84").Value = -34808
# Row 122
$ws.Range("H122").Value = 1980.1818
$ws.Range("I122").Value = 1980.1818
$ws.Range("K122").Value = 5940.5454
$ws.Range("M122").Value = -3490.5454
# Row 126
$ws.Range("H126").Value = 2165.125
$ws.Range("I126").Value = 2204.7144
$ws.Range("K126").Value = 6614.1432
$ws.Range("M126").Value = -4144.1432
# Row 132
$ws.Range("H132").Value = 33344626
$ws.Range("I132").Value = 55558332
$ws.Range("K132").Value = 166674996
$ws.Range("M132").Value = -166672466
